# Shift the delivery date/time values in columns C and D (rows 2-25) forward
# by 365 days (2021-08-18 -> 2022-08-18), keeping the time-of-day fraction
# unchanged, then leave the active selection on D4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$range = $ws.Range("C2:D25")

foreach ($cell in $range.Cells) {
    if ($cell.Value2 -ne $null) {
        $cell.Value = $cell.Value2 + 365
    }
}

$ws.Range("D4").Select()
